$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# Remove the animation effect that targets the soon-to-be-deleted "Group 161"
# shape, which also removes the now-orphaned <p:timing> block entirely.
if ($s.TimeLine.MainSequence.Count -gt 0) {
    $s.TimeLine.MainSequence.Item(1).Delete()
}

# Remove the whole "Group 161" shape (the small worked a-b-c-d-e-f example graph).
$s.Shapes.Item("Group 161").Delete()

# Update the numbers shown inside the remaining "Group 40" diagram.
$grp40 = $s.Shapes.Item("Group 40")

$tb = $grp40.GroupItems.Item("TextBox 67")
$tb.TextFrame.TextRange.Text = "4"
$tb.Width = 8.9617325

$grp40.GroupItems.Item("TextBox 70").TextFrame.TextRange.Text = "3"
$grp40.GroupItems.Item("TextBox 74").TextFrame.TextRange.Text = "2"
$grp40.GroupItems.Item("TextBox 75").TextFrame.TextRange.Text = "1"
$grp40.GroupItems.Item("TextBox 77").TextFrame.TextRange.Text = "0"
